$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(319).Insert()

$ws.Cells.Item(319, 1).Value = 10
$ws.Cells.Item(319, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(319, 3).Value = "La Araucanía"
$ws.Cells.Item(319, 4).Value = 44508
$ws.Cells.Item(319, 5).Value = 9
$ws.Cells.Item(319, 6).Value = 100112027
$ws.Cells.Item(319, 7).Value = "Melón"
$ws.Cells.Item(319, 8).Value = "Tuna"
$ws.Cells.Item(319, 9).Value = "Primera"
$ws.Cells.Item(319, 10).Value = 30
$ws.Cells.Item(319, 11).Value = 43000
$ws.Cells.Item(319, 12).Value = 43000
$ws.Cells.Item(319, 13).Value = 43000
$ws.Cells.Item(319, 14).Value = '$/caja 18 unidades'
$ws.Cells.Item(319, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(319, 16).Value = 2389
$ws.Cells.Item(319, 17).Value = 18
$ws.Cells.Item(319, 18).Value = "Hortaliza"
